$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in column C (Max column)
$ws.Range("C2").Value = 13.5
$ws.Range("C3").Value = 12
$ws.Range("C5").Value = 20

# Update the active selection to D3
$ws.Range("D3").Select()
